# Lernjournal angepasst und Zeit eingetragen
#
# The workbook originally holds two sheets:
#   "Lernjournal Aufgabe 2"  (first tab)
#   "Lernjournal Aufgabe 1"  (second tab)
#
# A new third task sheet ("Lernjournal Aufgabe 3") is inserted as the new
# first tab. It starts out as a copy of "Lernjournal Aufgabe 2" (same
# template/layout) and is then adapted: the logged activities/minutes/dates
# for the two work blocks are cleared out (ready for new entries) except for
# a single new entry that documents the layout/assignment rework, and the
# running totals are left to recompute accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert "Lernjournal Aufgabe 3" as a duplicate of "Lernjournal Aufgabe 2"
#    placed before it (i.e. as the new first sheet).
# ---------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("Lernjournal Aufgabe 2")
$sourceSheet.Copy($wb.Worksheets.Item(1))

$ws = $wb.Worksheets.Item(1)
$ws.Name = "Lernjournal Aufgabe 3"
$ws.Activate()

# ---------------------------------------------------------------------
# 2. Clear out the first block of logged entries (rows 11-15: activity
#    text, minutes spent, date) but keep the cell styling/borders intact.
# ---------------------------------------------------------------------
$block1 = $ws.Range("A11:C15")
foreach ($area in $block1.Areas) {
    $area.ClearContents()
}

# Rows 11-13 lose their custom row height (they no longer wrap multi-line
# text), rows 14-15 already used the default height.
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# 3. Clear out the second block of logged entries (rows 25-31), leaving
#    row 24 to be replaced with the new entry below.
# ---------------------------------------------------------------------
$block2 = $ws.Range("A25:C28,A29:C31")
foreach ($area in $block2.Areas) {
    $area.ClearContents()
}

$ws.Rows.Item(26).AutoFit()
$ws.Rows.Item(28).AutoFit()

# ---------------------------------------------------------------------
# 4. Replace row 24 with the new "layout reworked" entry.
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Dokument angepasst. Neue Aufgabenstellunge eingebaut und Lyout angepasst."
$ws.Range("B24").Value = 25
$ws.Range("C24").Value = 41433
$ws.Rows.Item(24).AutoFit()

# ---------------------------------------------------------------------
# 5. Update the sheet's remembered selection.
# ---------------------------------------------------------------------
$ws.Range("K48").Select()

Write-Output "Lernjournal Aufgabe 3 eingefuegt und Eintraege angepasst."
